$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column F (6th column), shifting old F (solution) to G
$ws.Range("F1").EntireColumn.Insert()

# New header for F1: "description"
$ws.Range("F1").Value = "description"

# New filter-question values for rows 2 and 3
$ws.Range("F2").Value = "coletor não funciona,não liga,não pega,parou de funcionar "
$ws.Range("F3").Value = "coletor travando,parado,travado,não funciona,não pega"

# Row heights for rows 2 and 3 match new content (54pt, autosized / no longer custom 39.75 / 18)
$ws.Rows.Item(2).RowHeight = 54
$ws.Rows.Item(3).RowHeight = 54

# New column F should share the same width as column E
$ws.Range("F1").ColumnWidth = $ws.Range("E1").ColumnWidth

# Update selection
$ws.Range("G6").Select()
